# Gatopolis Manager - "Desenvolvido nova tela de adicao de turmas"
#
# Applies:
#  - workbookView window geometry (xWindow/yWindow/windowWidth/windowHeight)
#  - shared-string casing fix: "PERIODO (Manhã/Tarde)" -> "PERIODO (MANHÃ/TARDE)"
#  - new active selection on the Manager sheet: A2 -> E2
#  - widened columns A:E and a brand-new column F (for the new "turma" field)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window geometry (bookViews/workbookView) ---------------------------
$win = $excel.ActiveWindow
$win.Left   = 0
$win.Top    = 0
$win.Width  = 25600
$win.Height = 16060

# --- Fix shared-string text (header for column E) ------------------------
$ws.Range("E1").Value = "PERIODO (MANHÃ/TARDE)"

# --- Column widths: widen existing columns + add new column F ------------
# (ColumnWidth is in "characters"; the stored OOXML <col width> attribute is
#  ColumnWidth + 5/6, so we back-solve for the desired stored widths.)
$ws.Columns.Item(1).ColumnWidth = 43.666666666666664   # -> stored width 44.5
$ws.Columns.Item(2).ColumnWidth = 36.330729166666664   # -> stored width ~37.1640625
$ws.Columns.Item(3).ColumnWidth = 22.166666666666668   # -> stored width 23
$ws.Columns.Item(4).ColumnWidth = 24.166666666666668   # -> stored width 25
$ws.Columns.Item(5).ColumnWidth = 22.998697916666668   # -> stored width ~23.83203125
$ws.Columns.Item(6).ColumnWidth = 15.998697916666666   # -> stored width ~16.83203125 (new col)

# --- Active cell selection moves from A2 to E2 ----------------------------
$ws.Range("E2").Select()
